# Add a new "Event sourcing" block right after the existing
# "Reify Persistence as Relationship (...)" paragraph, and before the
# blank paragraph that precedes "Sets Augmentation Domain:".
#
# New paragraphs (in order):
#   ""                                                                                  (blank)
#   "Event sourcing:"
#   ""                                                                                  (blank)
#   "(PersistenceeContext, PersistenceContext, PersistenceMember::new, PersistenceSubject);"
#   ""                                                                                  (blank)
#   "(PersistenceeContext, PersistenceContext, PersistenceMember::delete, PersistenceSubject);"
#   ""                                                                                  (blank)
#   "(PersistenceeContext, PersistenceSubject, PersistenceMember::delete, PersistenceSubject);"

$d = $word.ActiveDocument

$anchorText = "Reify Persistence as Relationship (Values as Relation Resources). Align domain / range with domains / primitive types (Member Kind, salary;ARS)."

$searchRange = $d.Content
$found = $searchRange.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Anchor paragraph not found"
}

# New content goes immediately after the anchor paragraph's text.
$cursor = $searchRange.End

$newLines = @(
    "",
    "Event sourcing:",
    "",
    "(PersistenceeContext, PersistenceContext, PersistenceMember::new, PersistenceSubject);",
    "",
    "(PersistenceeContext, PersistenceContext, PersistenceMember::delete, PersistenceSubject);",
    "",
    "(PersistenceeContext, PersistenceSubject, PersistenceMember::delete, PersistenceSubject);"
)

foreach ($line in $newLines) {
    $ip = $d.Range($cursor, $cursor)
    $ip.InsertParagraphAfter()
    $cursor = $cursor + 1
    if ($line -ne "") {
        $ip2 = $d.Range($cursor, $cursor)
        $ip2.InsertAfter($line)
        $cursor = $cursor + $line.Length
    }
}
